$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell used only as a formatting source (style index 3, never itself modified by this edit)
$donor = "C2"

$ws.Range("E2").Value = "2026-02-12 20:48:48"
$ws.Range("E3").Value = "2026-02-12 20:48:51"
$ws.Range("O3").Value = "-3.3 °C"
$ws.Range("E4").Value = "2026-02-12 20:48:53"
$ws.Range("J4").Value = "999.3 hPa"
$ws.Range("O4").Value = "16.2 °C"
$ws.Range("E5").Value = "2026-02-12 20:48:56"
$ws.Range("E6").Value = "2026-02-12 20:48:58"
$ws.Range("J6").Value = "999.1 hPa"
$ws.Range("K6").Value = "13.7 MJ/m2"
$ws.Range("O6").Value = "15.9 °C"
$ws.Range("E7").Value = "2026-02-12 20:49:01"
$ws.Range("J7").Value = "1001.9 hPa"
$ws.Range("E8").Value = "2026-02-12 20:49:03"
$ws.Range("J8").Value = "1001.3 hPa"
$ws.Range("E9").Value = "2026-02-12 20:49:06"
$ws.Range("H9").Value = "'66%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = "2026-02-12 20:49:09"
$ws.Range("O10").Value = "14.9 °C"
$ws.Range("E11").Value = "2026-02-12 20:49:11"
$ws.Range("O11").Value = "9.5 °C"
$ws.Range("E12").Value = "2026-02-12 20:49:14"
$ws.Range("H12").Value = "'71%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").Value = "2026-02-12 20:49:16"
$ws.Range("J13").Value = "1001.8 hPa"
$ws.Range("E14").Value = "2026-02-12 20:49:19"
$ws.Range("O14").Value = "17.0 °C"
$ws.Range("E15").Value = "2026-02-12 20:49:21"
$ws.Range("H15").Value = "'53%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = "2026-02-12 20:49:24"
$ws.Range("E17").Value = "2026-02-12 20:49:27"
$ws.Range("E18").Value = "2026-02-12 20:49:29"
$ws.Range("J18").Value = "999.6 hPa"
$ws.Range("O18").Value = "16.8 °C"
$ws.Range("E19").Value = "2026-02-12 20:49:32"
$ws.Range("O19").Value = "8.1 °C"
$ws.Range("E20").Value = "2026-02-12 20:49:34"
$ws.Range("H20").Value = "'85%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H20").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").Value = "2026-02-12 20:49:37"
$ws.Range("J21").Value = "1002.3 hPa"
$ws.Range("O21").Value = "9.3 °C"
$ws.Range("E22").Value = "2026-02-12 20:49:39"
$ws.Range("E23").Value = "2026-02-12 20:49:42"
$ws.Range("E24").Value = "2026-02-12 20:49:44"
$ws.Range("J24").Value = "1006.6 hPa"
$ws.Range("O24").Value = "11.6 °C"
$ws.Range("E25").Value = "2026-02-12 20:49:47"
$ws.Range("E26").Value = "2026-02-12 20:49:49"
$ws.Range("J26").Value = "998.7 hPa"
$ws.Range("E27").Value = "2026-02-12 20:49:52"
$ws.Range("E28").Value = "2026-02-12 20:49:55"
$ws.Range("J28").Value = "999.0 hPa"
$ws.Range("N28").Value = "7.6 °C 20:26 TU"
$ws.Range("O28").Value = "14.1 °C"
$ws.Range("E29").Value = "2026-02-12 20:49:57"
$ws.Range("H29").Value = "'58%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122) | Out-Null
$ws.Range("N29").Value = "6.7 °C 20:29 TU"
$ws.Range("O29").Value = "14.6 °C"
$ws.Range("E30").Value = "2026-02-12 20:50:00"
$ws.Range("J30").Value = "999.4 hPa"
$ws.Range("O30").Value = "12.2 °C"
$ws.Range("E31").Value = "2026-02-12 20:50:02"
$ws.Range("J31").Value = "998.8 hPa"
$ws.Range("E32").Value = "2026-02-12 20:50:05"
$ws.Range("E33").Value = "2026-02-12 20:50:07"
$ws.Range("J33").Value = "1001.5 hPa"
$ws.Range("O33").Value = "6.6 °C"
$ws.Range("E34").Value = "2026-02-12 20:50:10"
$ws.Range("E35").Value = "2026-02-12 20:50:13"
$ws.Range("J35").Value = "1008.1 hPa"
$ws.Range("E36").Value = "2026-02-12 20:50:15"
$ws.Range("H36").Value = "'60%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H36").PasteSpecial(-4122) | Out-Null
$ws.Range("J36").Value = "999.7 hPa"
$ws.Range("E37").Value = "2026-02-12 20:50:18"
$ws.Range("H37").Value = "'48%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H37").PasteSpecial(-4122) | Out-Null
$ws.Range("J37").Value = "1000.3 hPa"
$ws.Range("N37").Value = "5.7 °C 20:07 TU"
$ws.Range("O37").Value = "10.0 °C"
$ws.Range("E38").Value = "2026-02-12 20:50:20"
$ws.Range("E39").Value = "2026-02-12 20:50:22"
$ws.Range("E40").Value = "2026-02-12 20:50:25"
$ws.Range("H40").Value = "'55%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H40").PasteSpecial(-4122) | Out-Null
$ws.Range("J40").Value = "1003.1 hPa"
$ws.Range("O40").Value = "9.6 °C"
$ws.Range("E41").Value = "2026-02-12 20:50:27"
$ws.Range("J41").Value = "1005.6 hPa"
$ws.Range("E42").Value = "2026-02-12 20:50:30"
$ws.Range("H42").Value = "'61%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H42").PasteSpecial(-4122) | Out-Null
$ws.Range("N42").Value = "9.0 °C 20:23 TU"
$ws.Range("O42").Value = "14.2 °C"
$ws.Range("E43").Value = "2026-02-12 20:50:32"
$ws.Range("E44").Value = "2026-02-12 20:50:35"
$ws.Range("E45").Value = "2026-02-12 20:50:37"
$ws.Range("H45").Value = "'53%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H45").PasteSpecial(-4122) | Out-Null
$ws.Range("J45").Value = "1005.0 hPa"
$ws.Range("O45").Value = "7.0 °C"
$ws.Range("E46").Value = "2026-02-12 20:50:40"
$ws.Range("H46").Value = "'40%"
$ws.Range($donor).Copy() | Out-Null
$ws.Range("H46").PasteSpecial(-4122) | Out-Null
$ws.Range("J46").Value = "1007.4 hPa"
$ws.Range("K46").Value = "13.3 MJ/m2"
$ws.Range("O46").Value = "15.9 °C"
